$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "maa://24702 (94.22), maa://25390 (95.9), maa://36681 (87.01)"
$ws.Range("L2").Value = "*maa://24633 (55.48), *maa://30515 (69.0), *maa://34787 (71.21), ***maa://20792 (11.93), maa://39402 (87.18), ***maa://29083 (27.78)"
$ws.Range("T2").Value = "maa://22742 (91.45), *maa://20791 (62.86)"
$ws.Range("AB2").Value = "maa://21246 (91.32), maa://36684 (95.65), ***maa://22731 (6.67)"
$ws.Range("AF2").Value = "maa://25251 (92.39), ***maa://21730 (19.4), ***maa://39501 (15.0), *maa://36675 (60.0)"
$ws.Range("D3").Value = "maa://36987 (95.83), maa://40192 (100.0), maa://39849 (88.89)"
$ws.Range("L3").Value = "*maa://22880 (67.04), maa://20276 (84.71), *maa://22749 (66.67)"
$ws.Range("P3").Value = "maa://21249 (94.57), maa://26254 (95.83)"
$ws.Range("T3").Value = "maa://24617 (88.79), **maa://20790 (44.78), ***maa://37170 (18.52)"
$ws.Range("X3").Value = "maa://27396 (84.69), maa://27484 (96.04), maa://27480 (82.86)"
$ws.Range("T4").Value = "maa://32509 (97.92), maa://27295 (83.33), maa://22754 (91.67), *maa://21746 (55.81), *maa://31008 (78.05)"
$ws.Range("X4").Value = "**maa://32495 (47.88), ***maa://31785 (22.22), ***maa://36683 (28.26), maa://43217 (89.29)"
$ws.Range("AF4").Value = "*maa://30062 (61.36), ***maa://26209 (13.04), *maa://39394 (68.18)"
$ws.Range("L5").Value = "*maa://22757 (77.42)"
$ws.Range("D6").Value = "maa://42407 (93.94)"
$ws.Range("D7").Value = "maa://21955 (93.94)"
$ws.Range("L7").Value = "maa://28624 (93.02), maa://24957 (97.62)"
$ws.Range("T7").Value = "maa://21291 (85.71)"
$ws.Range("AF7").Value = "*maa://26191 (67.95), *maa://36671 (70.21), *maa://42530 (69.23)"
$ws.Range("A8").Value = "更新日期：2024.12.12 13:20:00"
$ws.Range("D8").Value = "*maa://21476 (71.74), **maa://39431 (45.45), *maa://37551 (57.14)"
$ws.Range("H8").Value = "*maa://24371 (54.29)"
$ws.Range("D9").Value = "maa://22765 (92.05), *maa://21915 (68.0)"
$ws.Range("L9").Value = "maa://22762 (91.86), maa://39552 (90.0)"
$ws.Range("P9").Value = "maa://22736 (81.11)"
$ws.Range("T9").Value = "**maa://22866 (30.19), maa://26222 (97.78)"
$ws.Range("AB9").Value = "maa://28711 (86.87), ***maa://22740 (5.88), **maa://39938 (48.0), **maa://27377 (42.86), ***maa://25174 (19.05), maa://40166 (93.33)"
$ws.Range("D10").Value = "***maa://25695 (19.44), **maa://32237 (40.48), ***maa://34206 (18.18), ***maa://39951 (17.5), ***maa://39243 (28.57)"
$ws.Range("P10").Value = "maa://28977 (91.36), maa://36669 (88.24), *maa://23264 (61.82)"
$ws.Range("T10").Value = "maa://27395 (95.93), maa://22755 (87.5), **maa://22756 (40.91), ***maa://21737 (10.61)"
$ws.Range("AF10").Value = "*maa://25021 (53.66), *maa://22733 (59.38), maa://22761 (100.0)"
$ws.Range("T11").Value = "maa://22747 (93.38), maa://22501 (98.53)"
$ws.Range("D12").Value = "maa://30766 (89.29), *maa://36678 (60.0)"
$ws.Range("X12").Value = "maa://22753 (91.25), *maa://21485 (76.87), maa://37962 (88.89)"
$ws.Range("AB12").Value = "maa://23669 (95.34), maa://36677 (93.75), maa://39872 (90.0)"
$ws.Range("D13").Value = "maa://24999 (91.73), maa://36673 (92.65), maa://25001 (85.51)"
$ws.Range("X13").Value = "*maa://34957 (79.03), *maa://22768 (51.61)"
$ws.Range("D14").Value = "maa://30764 (88.24)"
$ws.Range("L14").Value = "maa://26245 (96.43), maa://21288 (96.27), maa://39841 (94.87), maa://36682 (97.37)"
$ws.Range("P14").Value = "maa://23250 (98.62), maa://20107 (87.1), maa://22772 (100.0), **maa://22745 (50.0)"
$ws.Range("T14").Value = "maa://22521 (93.75), maa://42751 (100.0)"
$ws.Range("X14").Value = "maa://37468 (89.47)"
$ws.Range("AB14").Value = "maa://22764 (96.88)"
$ws.Range("D15").Value = "*maa://22743 (76.84), maa://22734 (83.9), *maa://30808 (64.52), ***maa://36048 (28.57)"
$ws.Range("H15").Value = "maa://24304 (88.56), maa://21478 (91.43)"
$ws.Range("AF15").Value = "maa://21364 (80.59), *maa://22766 (70.64), *maa://36666 (78.31)"
$ws.Range("D16").Value = "maa://21441 (96.28), maa://36679 (93.02), maa://37650 (96.88)"
$ws.Range("P16").Value = "maa://28504 (90.91)"
$ws.Range("T17").Value = "**maa://42324 (48.0)"
$ws.Range("H18").Value = "maa://24421 (90.25)"
$ws.Range("X18").Value = "maa://21917 (97.75), maa://22741 (83.33)"
$ws.Range("AB18").Value = "maa://24393 (97.56)"
$ws.Range("T19").Value = "maa://24386 (99.0)"
$ws.Range("AB19").Value = "*maa://30709 (62.78), *maa://36668 (54.67)"
$ws.Range("D20").Value = "maa://21432 (89.73), maa://25198 (93.0), *maa://20795 (51.18), maa://36680 (96.55)"
$ws.Range("X21").Value = "maa://20110 (86.76), maa://34946 (92.5)"
$ws.Range("AF21").Value = "maa://22524 (94.55), *maa://22432 (76.67)"
$ws.Range("L22").Value = "maa://27127 (85.15), *maa://22751 (73.85)"
$ws.Range("X22").Value = "maa://21282 (98.46), *maa://37649 (68.0)"
$ws.Range("D23").Value = "***maa://28036 (28.36), *maa://41753 (58.33)"
$ws.Range("L23").Value = "maa://39756 (93.78), maa://39875 (93.55)"
$ws.Range("P23").Value = "maa://30587 (91.89), *maa://29748 (75.59), ***maa://29785 (16.42), *maa://37566 (71.43)"
$ws.Range("AB23").Value = "maa://29652 (97.56)"
$ws.Range("X24").Value = "maa://29988 (86.61), maa://23504 (93.19), **maa://22892 (39.58), *maa://25141 (77.6), maa://36663 (81.82), ***maa://22815 (23.08)"
$ws.Range("AF24").Value = "maa://22523 (85.42), maa://36672 (80.77), maa://29910 (92.31), **maa://21440 (34.55)"
$ws.Range("H25").Value = "*maa://29063 (73.97), *maa://25311 (74.0), ***maa://22725 (4.84)"
$ws.Range("AB25").Value = "maa://31215 (85.42), *maa://24516 (79.55), maa://26001 (87.5)"
$ws.Range("D26").Value = "maa://41802 (91.67)"
$ws.Range("H26").Value = "maa://24913 (91.46)"
$ws.Range("AB26").Value = "maa://42235 (92.96)"
$ws.Range("AF26").Value = "maa://30511 (83.78), *maa://29760 (64.29)"
$ws.Range("T27").Value = "*maa://30624 (76.47)"
$ws.Range("AF27").Value = "maa://24023 (97.06)"
$ws.Range("D28").Value = "maa://24465 (90.76), maa://25725 (83.53)"
$ws.Range("T28").Value = "maa://23263 (94.85), *maa://29765 (60.53)"
$ws.Range("X28").Value = "maa://39929 (89.51), ***maa://39723 (14.29), maa://41749 (90.0)"
$ws.Range("L29").Value = "maa://28432 (93.11), *maa://28440 (76.6), maa://31400 (100.0), *maa://28650 (71.43)"
$ws.Range("P29").Value = "*maa://23168 (54.72), *maa://30050 (51.85)"
$ws.Range("AB30").Value = "maa://42979 (96.97)"
$ws.Range("L31").Value = "maa://35926 (93.92), maa://36258 (82.98), *maa://43904 (80.0)"
$ws.Range("T32").Value = "maa://42859 (95.52), maa://41108 (87.76), maa://41238 (96.1)"
$ws.Range("L35").Value = "maa://41296 (96.4)"
$ws.Range("AF38").Value = "maa://36697 (85.88)"
$ws.Range("P39").Value = "maa://24709 (92.0)"
$ws.Range("P41").Value = "**maa://35616 (38.24), maa://43177 (81.82)"
$ws.Range("H43").Value = "maa://22525 (92.59), maa://21284 (83.72)"
$ws.Range("H44").Value = "maa://29768 (97.79), maa://27728 (96.0)"
$ws.Range("H46").Value = "maa://35931 (92.42), maa://43901 (85.71)"
$ws.Range("H55").Value = "maa://32532 (92.22)"
$ws.Range("H57").Value = "maa://25176 (98.18)"
$ws.Range("H58").Value = "*maa://37964 (57.14)"
$ws.Range("H64").Value = "**maa://44405 (50.0)"
